$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.579.74"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.277.44"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "95.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -4.37%  "
$ws.Range("E10").Value = "  -8.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.02%  "
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "2.618.68"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.846"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "2.281.33"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "43.541.51"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000107"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.58%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "176.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("E33").Value = "  -3.33%  "
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.08%  "
$ws.Range("E40").Value = "  +8.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.236"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.91%  "
$ws.Range("E42").Value = "  +17.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("E45").Value = "  +4.32%  "
$ws.Range("E46").Value = "  -4.18%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("E50").Value = "  +7.93%  "
$ws.Range("D51").Value = "2.498.09"
$ws.Range("E51").Value = "  +1.51%  "
